$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 469.3
$ws.Range("I92").Value = 418.58823
$ws.Range("J92").Value = 756.6667
$ws.Range("K92").Value = 418.58823
$ws.Range("L92").Value = 756.6667
$ws.Range("M92").Value = 829.4117699999999
$ws.Range("N92").Value = -3252.6667
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("M94").Value = -2549
$ws.Range("H111").Value = 637.3889
$ws.Range("I111").Value = 542.9
$ws.Range("J111").Value = 755.5
$ws.Range("K111").Value = 1628.7
$ws.Range("L111").Value = 2266.5
$ws.Range("M111").Value = 1438.3
$ws.Range("N111").Value = -8400.5
$ws.Range("H132").Value = 5002877.5
$ws.Range("I132").Value = 6063181.5
$ws.Range("J132").Value = 4301.4287
$ws.Range("K132").Value = 18189544.5
$ws.Range("L132").Value = 12904.2861
$ws.Range("M132").Value = -18187014.5
$ws.Range("N132").Value = -17964.2861
$ws.Range("H137").Value = 3947.5642
$ws.Range("I137").Value = 3994.3
$ws.Range("J137").Value = 3791.7778
$ws.Range("K137").Value = 11982.9
$ws.Range("L137").Value = 11375.3334
$ws.Range("M137").Value = -9432.900000000001
$ws.Range("N137").Value = -16475.3334
$ws.Range("H138").Value = 2305.35
$ws.Range("I138").Value = 1390.4147
$ws.Range("J138").Value = 4279.684
$ws.Range("K138").Value = 4171.2441
$ws.Range("L138").Value = 12839.052
$ws.Range("M138").Value = 968.7559000000001
$ws.Range("N138").Value = -23119.052

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10463.948
$ws.Range("I32").Value = 9432.34
$ws.Range("K32").Value = 9432.34
$ws.Range("M32").Value = -9145.34
$ws.Range("H61").Value = 2547.8076
$ws.Range("I61").Value = 2207.25
$ws.Range("J61").Value = 2839.7144
$ws.Range("K61").Value = 2207.25
$ws.Range("L61").Value = 2839.7144
$ws.Range("M61").Value = -1995.25
$ws.Range("N61").Value = -3263.7144
$ws.Range("H74").Value = 2166.1
$ws.Range("I74").Value = 1488.8667
$ws.Range("J74").Value = 4197.8
$ws.Range("K74").Value = 1488.8667
$ws.Range("L74").Value = 4197.8
$ws.Range("M74").Value = -614.8667
$ws.Range("N74").Value = -5945.8
$ws.Range("H77").Value = 2166.1
$ws.Range("I77").Value = 1488.8667
$ws.Range("J77").Value = 4197.8
$ws.Range("K77").Value = 7444.333500000001
$ws.Range("L77").Value = 20989
$ws.Range("M77").Value = -3076.333500000001
$ws.Range("N77").Value = -29725
$ws.Range("H97").Value = 1599.8572
$ws.Range("I97").Value = 1533.1666
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1533.1666
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -1037.1666
$ws.Range("N97").Value = -2992
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1378
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 2947.9048
$ws.Range("I132").Value = 2493.375
$ws.Range("J132").Value = 4402.4
$ws.Range("K132").Value = 7480.125
$ws.Range("L132").Value = 13207.2
$ws.Range("M132").Value = -4950.125
$ws.Range("N132").Value = -18267.2
$ws.Range("H134").Value = 34900
$ws.Range("J134").Value = 34900
$ws.Range("L134").Value = 34900
$ws.Range("N134").Value = -45040
$ws.Range("H135").Value = 28441.389
$ws.Range("J135").Value = 28441.389
$ws.Range("L135").Value = 28441.389
$ws.Range("N135").Value = -38581.389
$ws.Range("H136").Value = 2547.8076
$ws.Range("I136").Value = 2207.25
$ws.Range("J136").Value = 2839.7144
$ws.Range("K136").Value = 6621.75
$ws.Range("L136").Value = 8519.143199999999
$ws.Range("M136").Value = -4071.75
$ws.Range("N136").Value = -13619.1432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 821.4
$ws.Range("I94").Value = 630.7143
$ws.Range("J94").Value = 1266.3334
$ws.Range("K94").Value = 630.7143
$ws.Range("L94").Value = 1266.3334
$ws.Range("M94").Value = -179.7143
$ws.Range("N94").Value = -2168.3334
$ws.Range("H99").Value = 4762.375
$ws.Range("I99").Value = 5519.8
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 5519.8
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -4021.8
$ws.Range("N99").Value = -6496
$ws.Range("H134").Value = 3060.875
$ws.Range("I134").Value = 3332.0344
$ws.Range("K134").Value = 9996.1032
$ws.Range("M134").Value = -7461.1032

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 55600
$ws.Range("J52").Value = 55600
$ws.Range("L52").Value = 55600
$ws.Range("N52").Value = -56188
$ws.Range("H58").Value = 14289713
$ws.Range("I58").Value = 3157.889
$ws.Range("J58").Value = 29416654
$ws.Range("K58").Value = 3157.889
$ws.Range("L58").Value = 29416654
$ws.Range("M58").Value = -2954.889
$ws.Range("N58").Value = -29417060
$ws.Range("H136").Value = 14289713
$ws.Range("I136").Value = 3157.889
$ws.Range("J136").Value = 29416654
$ws.Range("K136").Value = 9473.667000000001
$ws.Range("L136").Value = 88249962
$ws.Range("M136").Value = -6923.667000000001
$ws.Range("N136").Value = -88255062

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 833.6
$ws.Range("J18").Value = 1665.3334
$ws.Range("L18").Value = 4996.0002
$ws.Range("N18").Value = -5334.0002
$ws.Range("H131").Value = 1346.6394
$ws.Range("J131").Value = 1123.1
$ws.Range("L131").Value = 3369.3
$ws.Range("N131").Value = -13449.3

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 28626.25
$ws.Range("J19").Value = 36668.332
$ws.Range("L19").Value = 36668.332
$ws.Range("N19").Value = -37244.332
$ws.Range("H97").Value = 2407.7273
$ws.Range("I97").Value = 2348.5
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 2348.5
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1852.5
$ws.Range("N97").Value = -3992
$ws.Range("H122").Value = 7637.1816
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 6751.125
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 20253.375
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = -25153.375
$ws.Range("H124").Value = 31571.428
$ws.Range("J124").Value = 31571.428
$ws.Range("L124").Value = 31571.428
$ws.Range("N124").Value = -41391.428
$ws.Range("H132").Value = 3934.25
$ws.Range("I132").Value = 4708.1333
$ws.Range("J132").Value = 3381.476
$ws.Range("K132").Value = 14124.3999
$ws.Range("L132").Value = 10144.428
$ws.Range("M132").Value = -11594.3999
$ws.Range("N132").Value = -15204.428
$ws.Range("H134").Value = 31546.857
$ws.Range("J134").Value = 31546.857
$ws.Range("L134").Value = 94640.571
$ws.Range("N134").Value = -99710.571

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 372145.72
$ws.Range("I14").Value = 2500000
$ws.Range("K14").Value = 2500000
$ws.Range("M14").Value = -2499828
$ws.Range("H61").Value = 2247.7273
$ws.Range("I61").Value = 580.55554
$ws.Range("K61").Value = 580.55554
$ws.Range("M61").Value = -378.55554
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 1763.2667
$ws.Range("I100").Value = 1677.1818
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1677.1818
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1136.1818
$ws.Range("N100").Value = -3082
$ws.Range("H113").Value = 2247.7273
$ws.Range("I113").Value = 580.55554
$ws.Range("K113").Value = 580.55554
$ws.Range("M113").Value = 1589.44446
$ws.Range("H136").Value = 3784
$ws.Range("I136").Value = 4009.3333
$ws.Range("K136").Value = 12027.9999
$ws.Range("M136").Value = -9477.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 85065
$ws.Range("J48").Value = 85065
$ws.Range("L48").Value = 85065
$ws.Range("N48").Value = -86203
$ws.Range("H100").Value = 627
$ws.Range("I100").Value = 536
$ws.Range("J100").Value = 718
$ws.Range("K100").Value = 1072
$ws.Range("L100").Value = 1436
$ws.Range("M100").Value = -531
$ws.Range("N100").Value = -2518
$ws.Range("H132").Value = 5760.0557
$ws.Range("I132").Value = 2517.8462
$ws.Range("J132").Value = 14189.8
$ws.Range("K132").Value = 7553.5386
$ws.Range("L132").Value = 42569.39999999999
$ws.Range("M132").Value = -5023.5386
$ws.Range("N132").Value = -47629.39999999999
$ws.Range("H136").Value = 3358.3809
$ws.Range("I136").Value = 3547.6365
$ws.Range("J136").Value = 3150.2
$ws.Range("K136").Value = 10642.9095
$ws.Range("L136").Value = 9450.599999999999
$ws.Range("M136").Value = -8092.9095
$ws.Range("N136").Value = -14550.6

